$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (Volume/Number and date range) ---
$ws.Cells.Item(8, 1).Value = "Volume 31   Number  46"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  11/11/2024  Through  11/17/2024"

# --- Column width adjustments for I (9) and J (10) to match H (8) ---
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth()
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(8).ColumnWidth()

# --- String-type cell changes (copy style+value from an existing "0"/"***.*" cell) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F29").Copy($ws.Range("G29"))
$ws.Range("E29").Copy($ws.Range("H29"))
$ws.Range("F30").Copy($ws.Range("G30"))
$ws.Range("E30").Copy($ws.Range("H30"))

# --- Numeric cell value changes ---
$ws.Range("L14").Value = -75
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 38
$ws.Range("K15").Value = 15.151515151515
$ws.Range("L15").Value = 11.764705882352
$ws.Range("M15").Value = 65.217391304347
$ws.Range("N15").Value = 26.666666666666
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 56
$ws.Range("G16").Value = 38
$ws.Range("H16").Value = 47.368421052631
$ws.Range("I16").Value = 484
$ws.Range("J16").Value = 413
$ws.Range("K16").Value = 17.191283292978
$ws.Range("L16").Value = 34.444444444444
$ws.Range("M16").Value = 49.84520123839
$ws.Range("N16").Value = -64.490095377843
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 56
$ws.Range("H17").Value = -17.857142857142
$ws.Range("I17").Value = 721
$ws.Range("J17").Value = 657
$ws.Range("K17").Value = 9.741248097412
$ws.Range("L17").Value = 51.470588235294
$ws.Range("M17").Value = 227.727272727273
$ws.Range("N17").Value = 65.366972477064
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 16
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 214
$ws.Range("J18").Value = 168
$ws.Range("K18").Value = 27.380952380952
$ws.Range("L18").Value = 48.611111111111
$ws.Range("M18").Value = -15.748031496063
$ws.Range("N18").Value = -88.923395445134
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 80
$ws.Range("H19").Value = -6.25
$ws.Range("I19").Value = 983
$ws.Range("J19").Value = 943
$ws.Range("K19").Value = 4.24178154825
$ws.Range("L19").Value = -6.291706387035
$ws.Range("M19").Value = 103.099173553719
$ws.Range("N19").Value = -10.47358834244
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -44.444444444444
$ws.Range("F20").Value = 19
$ws.Range("H20").Value = -26.923076923076
$ws.Range("I20").Value = 256
$ws.Range("J20").Value = 303
$ws.Range("K20").Value = -15.511551155115
$ws.Range("L20").Value = 28
$ws.Range("M20").Value = 96.923076923076
$ws.Range("N20").Value = -87.238285144566
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 51
$ws.Range("E21").Value = 5.882352941176
$ws.Range("F21").Value = 217
$ws.Range("G21").Value = 218
$ws.Range("H21").Value = -0.45871559633
$ws.Range("I21").Value = 2698
$ws.Range("J21").Value = 2520
$ws.Range("K21").Value = 7.063492063492
$ws.Range("L21").Value = 18.802289740202
$ws.Range("M21").Value = 87.752261656228
$ws.Range("N21").Value = -60.836115546523
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 55
$ws.Range("K22").Value = -20
$ws.Range("M22").Value = 62.962962962963
$ws.Range("C24").Value = 48
$ws.Range("D24").Value = 67
$ws.Range("E24").Value = -28.358208955223
$ws.Range("G24").Value = 210
$ws.Range("H24").Value = -18.571428571428
$ws.Range("I24").Value = 2644
$ws.Range("J24").Value = 2488
$ws.Range("K24").Value = 6.270096463022
$ws.Range("L24").Value = 26.386233269598
$ws.Range("M24").Value = 65.872020075282
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = -24.489795918367
$ws.Range("F25").Value = 120
$ws.Range("G25").Value = 162
$ws.Range("H25").Value = -25.925925925925
$ws.Range("I25").Value = 2054
$ws.Range("J25").Value = 1839
$ws.Range("K25").Value = 11.691136487221
$ws.Range("L25").Value = 37.024683122081
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 109
$ws.Range("G26").Value = 99
$ws.Range("H26").Value = 10.10101010101
$ws.Range("I26").Value = 1224
$ws.Range("J26").Value = 1027
$ws.Range("K26").Value = 19.182083739045
$ws.Range("L26").Value = 53.575909661229
$ws.Range("M26").Value = 117.021276595745
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 14.285714285714
$ws.Range("I27").Value = 58
$ws.Range("K27").Value = 11.538461538461
$ws.Range("L27").Value = 9.43396226415
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 225
$ws.Range("I28").Value = 140
$ws.Range("K28").Value = 18.64406779661
$ws.Range("L28").Value = 20.689655172413
$ws.Range("L29").Value = -85.714285714285
$ws.Range("N29").Value = -96.428571428571
$ws.Range("L30").Value = -77.777777777777
$ws.Range("N30").Value = -95.918367346938
